$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version: 2.1.0 -> 2.2.0-ballot
$wsMeta.Range("B3").Value = "2.2.0-ballot"

# Date: 2025-12-19T08:22:07+00:00 -> 2025-12-19T09:47:21+00:00
$wsMeta.Range("B8").Value = "2025-12-19T09:47:21+00:00"

# Base Definition: append |4.0.1
$wsMeta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension|4.0.1"

# --- Elements sheet updates ---
$wsElem = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s): append |2.2.0-ballot to the referenced profile URL,
# keeping the trailing newline that was already present in the cell text.
$wsElem.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-document-reference|2.2.0-ballot)`n"

# Column K width grew (Excel "best fit" recalculation) as a side effect of the
# longer text now in K6. Set it to the closest width this engine can produce.
$wsElem.Columns.Item(11).ColumnWidth = 88.8
